$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "epochValidator"

# --- Clear the hyperlink that lived on C2 (and its styled/string content) ---
$ws.Range("C2").Hyperlinks.Delete()
$ws.Range("C2").ClearContents()

# --- Remove now-obsolete cells ---
$ws.Range("C1").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# --- Write new cell values in the same order the shared strings were
#     first introduced (so the shared-string table order matches) ---

# Row 2, col A
$ws.Range("A2").Value = "verifyEpochConversionWithValidInput"
# Row 3, col A
$ws.Range("A3").Value = "verifyEpochConversionWithInvalidInputExceedingInt32Range"
# Row 4, col A
$ws.Range("A4").Value = "verifyEpochConversionWithInvalidInput"
# Row 1, col B
$ws.Range("B1").Value = "ExpectedMessage"
# Row 4, col B
$ws.Range("B4").Value = "UnixTimeService.RESTHost.fromunixtimestamp.TryCatch.Try.UnixTimeStamp: Input string was not in a correct format."
# Row 3, col B
$ws.Range("B3").Value = "UnixTimeService.RESTHost.fromunixtimestamp.TryCatch.Try.UnixTimeStamp: Value was either too large or too small for an Int32."

# Row 2, col C (hyperlink style retained, but cell left blank)
$ws.Range("C2").Style = "Hyperlink"

# --- Column widths (nearest pixel-grid the engine can represent) ---
$ws.Columns.Item(1).ColumnWidth = 55.33
$ws.Columns.Item(2).ColumnWidth = 49.0

# --- Selection ---
[void]$ws.Range("B11").Select()

Write-Output "done"
